# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The four detail rows (16-19) on "Hoja1" are reshuffled into a new order
# (same combinations of worker / period, just re-sequenced) and the
# "Salario Basico" (column G) is refreshed from 1,000,000 to 908,526 for
# every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New values for the detail block B16:G19 (Tipo Doc, N Doc, Nombre, Periodo, Valor Mora, Salario Basico)
$data = @(
    @("CC", "1063144788", "SAMIRA DEL CARMEN DIAZ VARGAS", "2208", 36341, 908526),
    @("CC", "1062674021", "GREGORIO JOSE LUNA FLOREZ",     "2208", 36341, 908526),
    @("CC", "1063144788", "SAMIRA DEL CARMEN DIAZ VARGAS", "2209", 16959, 908526),
    @("CC", "1062674021", "GREGORIO JOSE LUNA FLOREZ",     "2209", 16959, 908526)
)

$row = 16
foreach ($line in $data) {
    $ws.Cells.Item($row, 2).Value = $line[0]   # B - Tipo Doc Trabajador
    $ws.Cells.Item($row, 3).Value = $line[1]   # C - N Doc Trabajador
    $ws.Cells.Item($row, 4).Value = $line[2]   # D - Nombre Trabajador
    $ws.Cells.Item($row, 5).Value = $line[3]   # E - Periodo Mora
    $ws.Cells.Item($row, 6).Value = $line[4]   # F - Valor Mora
    $ws.Cells.Item($row, 7).Value = $line[5]   # G - Salario Basico
    $row++
}
